$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.366.29'
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').Value = '1.825.15'
$ws.Range('E3').Value = '  +2.76%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '317.18'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5346'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4044'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +8.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07609'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.45%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.85'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('E11').Value = '  +1.60%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.324'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +4.38%  '
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.615'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +5.80%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '20.79'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.04%  '
$ws.Range('D16').Value = '1.824.19'
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '89.27'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.79%  '
$ws.Range('E18').Value = '  +2.21%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06604'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.21%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.64'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.53%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.111'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.14%  '
$ws.Range('D23').Value = '28.386.22'
$ws.Range('E23').Value = '  +1.53%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.17'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.97%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.192'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +5.65%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.460'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +8.36%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '157.72'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.10%  '
$ws.Range('E28').Value = '  +1.85%  '
$ws.Range('D29').Value = '2.035.81'
$ws.Range('E29').Value = '  +2.71%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '123.70'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +3.45%  '
$ws.Range('E31').Value = '  +1.22%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1097'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +5.34%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.651'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.82%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07277'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +14.65%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.2234'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02342'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.66%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '8.865'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +5.90%  '
$ws.Range('E39').Value = '  +4.91%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6247'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.36%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '11.26'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.84%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.181'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.77%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.000'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.400'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.66%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.52'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.08%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.705'
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5837'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.95%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '125.19'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.988'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.68%  '
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06888'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.49%  '
